# "Set up search male thread."
# Populate the question list: keep the original question (reworded with a
# trailing "?") and add the new O/P-drive questions around it, then drop the
# now-unused A1 row so the data starts at row 2 and the sheet's used range
# becomes A2:A8. Finally widen column A to fit the longer text and leave the
# selection on the next empty row (A9), matching the post-edit sheet view.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A3").Value = "O drive max = 255 (max pixel brightness) ?"
$ws.Range("A2").Value = "The male/neopixel/ring always white ?"
$ws.Range("A4").Value = "O drive max = 0 (pixel shutdown) ?"
$ws.Range("A5").Value = "P drive max = 255 (max pixel brightness) ?"
$ws.Range("A6").Value = "P drive max = 0 (pixel shutdown) ?"
$ws.Range("A7").Value = "P drive increment X unit every Y seconds. X and Y ?"
$ws.Range("A8").Value = "O drive increment X unit every Y seconds. X and Y ?"

# The old A1 cell is no longer needed now that the text lives in A2.
$ws.Range("A1").Value = ""

# Widen column A (no longer a "best fit" autofit column) to fit the new text.
$ws.Columns("A").ColumnWidth = 47.666666666666664

# Leave the selection on the next free row, like the author did after typing.
$ws.Range("A9").Select()
